$d = $word.ActiveDocument

# Locate the "MIC" book-code heading (Heading2) that precedes the
# italic "Michée" paragraph we need to remove.
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
$rng = $d.Content
$found = $rng.Find.Execute("MIC", $true, $true, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'MIC' heading paragraph"
}

$micIndex = $rng.Paragraphs.Item(1).Index
$micPara = $d.Paragraphs.Item($micIndex)

# The paragraph right after "MIC" is the standalone italic "Michée"
# paragraph; delete it in full (text + its paragraph mark) so the
# following paragraph merges up, exactly like the authored edit.
$michaeePara = $micPara.Next()
$michaeePara.Range.Delete() | Out-Null
